$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.600915312767029
$ws.Range("B1").Value = 4.229971885681152
$ws.Range("C1").Value = 3.899255037307739
$ws.Range("D1").Value = 1.818310022354126
$ws.Range("E1").Value = 1.053126573562622
